$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-22 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-23 Sunday", 2) | Out-Null
$d.Content.Find.Execute("36+63=", $true, $false, $false, $false, $false, $true, 1, $false, "41-37=", 2) | Out-Null
$d.Content.Find.Execute("38+54=", $true, $false, $false, $false, $false, $true, 1, $false, "40+47=", 2) | Out-Null
$d.Content.Find.Execute("60-5=", $true, $false, $false, $false, $false, $true, 1, $false, "5+3=", 2) | Out-Null
$d.Content.Find.Execute("94-19=", $true, $false, $false, $false, $false, $true, 1, $false, "14+51=", 2) | Out-Null
$d.Content.Find.Execute("93-8=", $true, $false, $false, $false, $false, $true, 1, $false, "14+55=", 2) | Out-Null
$d.Content.Find.Execute("86-59=", $true, $false, $false, $false, $false, $true, 1, $false, "16-14=", 2) | Out-Null
$d.Content.Find.Execute("73-40=", $true, $false, $false, $false, $false, $true, 1, $false, "17+60=", 2) | Out-Null
$d.Content.Find.Execute("78+4=", $true, $false, $false, $false, $false, $true, 1, $false, "52-39=", 2) | Out-Null
$d.Content.Find.Execute("84-32=", $true, $false, $false, $false, $false, $true, 1, $false, "64+17=", 2) | Out-Null
$d.Content.Find.Execute("28+22=", $true, $false, $false, $false, $false, $true, 1, $false, "29-23=", 2) | Out-Null
$d.Content.Find.Execute("14-6=", $true, $false, $false, $false, $false, $true, 1, $false, "54-36=", 2) | Out-Null
$d.Content.Find.Execute("20+17=", $true, $false, $false, $false, $false, $true, 1, $false, "42+10=", 2) | Out-Null
$d.Content.Find.Execute("24+45=", $true, $false, $false, $false, $false, $true, 1, $false, "55+29=", 2) | Out-Null
$d.Content.Find.Execute("67+30=", $true, $false, $false, $false, $false, $true, 1, $false, "45-17=", 2) | Out-Null
$d.Content.Find.Execute("53-49=", $true, $false, $false, $false, $false, $true, 1, $false, "78-22=", 2) | Out-Null
$d.Content.Find.Execute("26+56=", $true, $false, $false, $false, $false, $true, 1, $false, "90-52=", 2) | Out-Null
$d.Content.Find.Execute("17-14=", $true, $false, $false, $false, $false, $true, 1, $false, "61+23=", 2) | Out-Null
$d.Content.Find.Execute("21+24=", $true, $false, $false, $false, $false, $true, 1, $false, "17+66=", 2) | Out-Null
$d.Content.Find.Execute("76-62=", $true, $false, $false, $false, $false, $true, 1, $false, "76+3=", 2) | Out-Null
$d.Content.Find.Execute("6+70=", $true, $false, $false, $false, $false, $true, 1, $false, "83-47=", 2) | Out-Null
$d.Content.Find.Execute("65+0=", $true, $false, $false, $false, $false, $true, 1, $false, "22+39=", 2) | Out-Null
$d.Content.Find.Execute("59-44=", $true, $false, $false, $false, $false, $true, 1, $false, "53-7=", 2) | Out-Null
$d.Content.Find.Execute("87-75=", $true, $false, $false, $false, $false, $true, 1, $false, "37-30=", 2) | Out-Null
$d.Content.Find.Execute("40-28=", $true, $false, $false, $false, $false, $true, 1, $false, "98-93=", 2) | Out-Null
$d.Content.Find.Execute("64-3=", $true, $false, $false, $false, $false, $true, 1, $false, "25+63=", 2) | Out-Null
$d.Content.Find.Execute("23-23=", $true, $false, $false, $false, $false, $true, 1, $false, "88-57=", 2) | Out-Null
$d.Content.Find.Execute("10+33=", $true, $false, $false, $false, $false, $true, 1, $false, "87-43=", 2) | Out-Null
$d.Content.Find.Execute("59-32=", $true, $false, $false, $false, $false, $true, 1, $false, "62+37=", 2) | Out-Null
$d.Content.Find.Execute("25+3=", $true, $false, $false, $false, $false, $true, 1, $false, "81-28=", 2) | Out-Null
$d.Content.Find.Execute("20-12=", $true, $false, $false, $false, $false, $true, 1, $false, "12+70=", 2) | Out-Null
$d.Content.Find.Execute("52+18=", $true, $false, $false, $false, $false, $true, 1, $false, "38-19=", 2) | Out-Null
$d.Content.Find.Execute("35-22=", $true, $false, $false, $false, $false, $true, 1, $false, "27+15=", 2) | Out-Null
$d.Content.Find.Execute("32-9=", $true, $false, $false, $false, $false, $true, 1, $false, "80-26=", 2) | Out-Null
$d.Content.Find.Execute("17+18=", $true, $false, $false, $false, $false, $true, 1, $false, "57+0=", 2) | Out-Null
$d.Content.Find.Execute("27+5=", $true, $false, $false, $false, $false, $true, 1, $false, "23-4=", 2) | Out-Null
$d.Content.Find.Execute("94-69=", $true, $false, $false, $false, $false, $true, 1, $false, "11+10=", 2) | Out-Null
$d.Content.Find.Execute("75-60=", $true, $false, $false, $false, $false, $true, 1, $false, "95-5=", 2) | Out-Null
$d.Content.Find.Execute("49-44=", $true, $false, $false, $false, $false, $true, 1, $false, "43-27=", 2) | Out-Null
$d.Content.Find.Execute("38+55=", $true, $false, $false, $false, $false, $true, 1, $false, "24+56=", 2) | Out-Null
$d.Content.Find.Execute("49+5=", $true, $false, $false, $false, $false, $true, 1, $false, "21+61=", 2) | Out-Null
$d.Content.Find.Execute("60-41=", $true, $false, $false, $false, $false, $true, 1, $false, "93-76=", 2) | Out-Null
$d.Content.Find.Execute("18+23=", $true, $false, $false, $false, $false, $true, 1, $false, "34-13=", 2) | Out-Null
$d.Content.Find.Execute("75-40=", $true, $false, $false, $false, $false, $true, 1, $false, "65-49=", 2) | Out-Null
$d.Content.Find.Execute("96-31=", $true, $false, $false, $false, $false, $true, 1, $false, "87+2=", 2) | Out-Null
$d.Content.Find.Execute("1+38=", $true, $false, $false, $false, $false, $true, 1, $false, "63-13=", 2) | Out-Null
$d.Content.Find.Execute("4+64=", $true, $false, $false, $false, $false, $true, 1, $false, "68-56=", 2) | Out-Null
$d.Content.Find.Execute("34-5=", $true, $false, $false, $false, $false, $true, 1, $false, "88-72=", 2) | Out-Null
$d.Content.Find.Execute("82+2=", $true, $false, $false, $false, $false, $true, 1, $false, "92-71=", 2) | Out-Null
$d.Content.Find.Execute("36+15=", $true, $false, $false, $false, $false, $true, 1, $false, "1+45=", 2) | Out-Null
$d.Content.Find.Execute("25+67=", $true, $false, $false, $false, $false, $true, 1, $false, "7+20=", 2) | Out-Null
$d.Content.Find.Execute("77-32=", $true, $false, $false, $false, $false, $true, 1, $false, "68+14=", 2) | Out-Null
$d.Content.Find.Execute("6+42=", $true, $false, $false, $false, $false, $true, 1, $false, "13+70=", 2) | Out-Null
$d.Content.Find.Execute("32+35=", $true, $false, $false, $false, $false, $true, 1, $false, "42-13=", 2) | Out-Null
$d.Content.Find.Execute("34+20=", $true, $false, $false, $false, $false, $true, 1, $false, "54+35=", 2) | Out-Null
$d.Content.Find.Execute("27-19=", $true, $false, $false, $false, $false, $true, 1, $false, "17+5=", 2) | Out-Null
$d.Content.Find.Execute("99-84=", $true, $false, $false, $false, $false, $true, 1, $false, "87+2=", 2) | Out-Null
$d.Content.Find.Execute("67+23=", $true, $false, $false, $false, $false, $true, 1, $false, "83-80=", 2) | Out-Null
$d.Content.Find.Execute("7+23=", $true, $false, $false, $false, $false, $true, 1, $false, "84-38=", 2) | Out-Null
$d.Content.Find.Execute("2+96=", $true, $false, $false, $false, $false, $true, 1, $false, "27+27=", 2) | Out-Null
$d.Content.Find.Execute("19-1=", $true, $false, $false, $false, $false, $true, 1, $false, "5+25=", 2) | Out-Null
$d.Content.Find.Execute("52-33=", $true, $false, $false, $false, $false, $true, 1, $false, "97-42=", 2) | Out-Null
$d.Content.Find.Execute("2+26=", $true, $false, $false, $false, $false, $true, 1, $false, "58-45=", 2) | Out-Null
$d.Content.Find.Execute("52+4=", $true, $false, $false, $false, $false, $true, 1, $false, "91-79=", 2) | Out-Null
$d.Content.Find.Execute("43+22=", $true, $false, $false, $false, $false, $true, 1, $false, "52-21=", 2) | Out-Null
$d.Content.Find.Execute("50+12=", $true, $false, $false, $false, $false, $true, 1, $false, "55+12=", 2) | Out-Null
$d.Content.Find.Execute("9+24=", $true, $false, $false, $false, $false, $true, 1, $false, "86+4=", 2) | Out-Null
$d.Content.Find.Execute("24+24=", $true, $false, $false, $false, $false, $true, 1, $false, "70-8=", 2) | Out-Null
$d.Content.Find.Execute("41+26=", $true, $false, $false, $false, $false, $true, 1, $false, "7+33=", 2) | Out-Null
$d.Content.Find.Execute("7+91=", $true, $false, $false, $false, $false, $true, 1, $false, "57+37=", 2) | Out-Null
$d.Content.Find.Execute("48-39=", $true, $false, $false, $false, $false, $true, 1, $false, "28+9=", 2) | Out-Null
$d.Content.Find.Execute("73+0=", $true, $false, $false, $false, $false, $true, 1, $false, "89-63=", 2) | Out-Null
$d.Content.Find.Execute("64-41=", $true, $false, $false, $false, $false, $true, 1, $false, "86-58=", 2) | Out-Null
$d.Content.Find.Execute("96-77=", $true, $false, $false, $false, $false, $true, 1, $false, "13+40=", 2) | Out-Null
$d.Content.Find.Execute("45-6=", $true, $false, $false, $false, $false, $true, 1, $false, "72-22=", 2) | Out-Null
$d.Content.Find.Execute("9+90=", $true, $false, $false, $false, $false, $true, 1, $false, "22-4=", 2) | Out-Null
$d.Content.Find.Execute("54+19=", $true, $false, $false, $false, $false, $true, 1, $false, "50-35=", 2) | Out-Null
$d.Content.Find.Execute("48+10=", $true, $false, $false, $false, $false, $true, 1, $false, "21+75=", 2) | Out-Null
$d.Content.Find.Execute("35+1=", $true, $false, $false, $false, $false, $true, 1, $false, "92-47=", 2) | Out-Null
$d.Content.Find.Execute("64-61=", $true, $false, $false, $false, $false, $true, 1, $false, "63-0=", 2) | Out-Null
$d.Content.Find.Execute("59+30=", $true, $false, $false, $false, $false, $true, 1, $false, "10+67=", 2) | Out-Null
$d.Content.Find.Execute("69-55=", $true, $false, $false, $false, $false, $true, 1, $false, "74+14=", 2) | Out-Null
$d.Content.Find.Execute("95-76=", $true, $false, $false, $false, $false, $true, 1, $false, "90-85=", 2) | Out-Null
$d.Content.Find.Execute("74+24=", $true, $false, $false, $false, $false, $true, 1, $false, "62-16=", 2) | Out-Null
$d.Content.Find.Execute("64-55=", $true, $false, $false, $false, $false, $true, 1, $false, "46+42=", 2) | Out-Null
$d.Content.Find.Execute("47-22=", $true, $false, $false, $false, $false, $true, 1, $false, "4+5=", 2) | Out-Null
$d.Content.Find.Execute("95-10=", $true, $false, $false, $false, $false, $true, 1, $false, "26+55=", 2) | Out-Null
$d.Content.Find.Execute("95+3=", $true, $false, $false, $false, $false, $true, 1, $false, "77-34=", 2) | Out-Null
$d.Content.Find.Execute("72-70=", $true, $false, $false, $false, $false, $true, 1, $false, "79-77=", 2) | Out-Null
$d.Content.Find.Execute("10+27=", $true, $false, $false, $false, $false, $true, 1, $false, "73-52=", 2) | Out-Null
$d.Content.Find.Execute("46-22=", $true, $false, $false, $false, $false, $true, 1, $false, "4+66=", 2) | Out-Null
$d.Content.Find.Execute("61-31=", $true, $false, $false, $false, $false, $true, 1, $false, "6+14=", 2) | Out-Null
$d.Content.Find.Execute("42-21=", $true, $false, $false, $false, $false, $true, 1, $false, "12+32=", 2) | Out-Null
$d.Content.Find.Execute("2+57=", $true, $false, $false, $false, $false, $true, 1, $false, "88-77=", 2) | Out-Null
$d.Content.Find.Execute("53+43=", $true, $false, $false, $false, $false, $true, 1, $false, "55-4=", 2) | Out-Null
$d.Content.Find.Execute("78-18=", $true, $false, $false, $false, $false, $true, 1, $false, "18+12=", 2) | Out-Null
$d.Content.Find.Execute("37-7=", $true, $false, $false, $false, $false, $true, 1, $false, "70+17=", 2) | Out-Null
$d.Content.Find.Execute("6-2=", $true, $false, $false, $false, $false, $true, 1, $false, "34+49=", 2) | Out-Null
$d.Content.Find.Execute("15+52=", $true, $false, $false, $false, $false, $true, 1, $false, "34+10=", 2) | Out-Null
$d.Content.Find.Execute("83-77=", $true, $false, $false, $false, $false, $true, 1, $false, "48-48=", 2) | Out-Null
$d.Content.Find.Execute("68+10=", $true, $false, $false, $false, $false, $true, 1, $false, "92-11=", 2) | Out-Null
